$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEC-2020")

# Rows 6:8 form a template (merged C6:G7 "blank" block + a filled data row 8).
# Copy that block down into rows 13:15 to pick up the exact same cell styles
# (including the C13:G14 merge) that the real edit re-used.
$ws.Range("A6:G8").Copy($ws.Range("A13"))

# Restore the distinctive row heights used by the merged block.
$ws.Rows("13").RowHeight = 25.95
$ws.Rows("14").RowHeight = 21

# Fill in the new day entries.
$ws.Range("B13").Value = 44177
$ws.Range("B14").Value = 44178

$ws.Range("B15").Value = "2020/12/014"
$ws.Range("C15").Value = "QMVAR 2.0"
$ws.Range("D15").Value = "Design issue fixing"

# F15 needs the "WIP" styling (blue fill), so pull formatting from another
# WIP cell before writing the text in.
$ws.Range("F12").Copy($ws.Range("F15"))
$ws.Range("F15").Value = "WIP"

# Update the active selection to match the saved view state.
$ws.Range("C20").Select()
